$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column: force text by leading apostrophe, then reset style so no
# quotePrefix/format bleeds into the cell (matches source: plain inlineStr,
# no "s" attribute).
# E column: already non-numeric (leading/trailing spaces), so plain Value works.

$ws.Range("D2").Value = "'" + '76.138.85'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = "'" + '2.919.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.53%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'" + '199.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.18%  '
$ws.Range("D6").Value = "'" + '600.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = "'" + '0.551'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.04%  '
$ws.Range("E9").Value = '  +4.59%  '
$ws.Range("D10").Value = "'" + '2.918.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.59%  '
$ws.Range("D11").Value = "'" + '0.431'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +16.95%  '
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = "'" + '4.89'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = "'" + '3.455.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.60%  '
$ws.Range("D15").Value = "'" + '76.006.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.56%  '
$ws.Range("D16").Value = "'" + '0.0000192'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.94%  '
$ws.Range("D17").Value = "'" + '27.56'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.26%  '
$ws.Range("D18").Value = "'" + '2.910.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.46%  '
$ws.Range("D19").Value = "'" + '13.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.99%  '
$ws.Range("D20").Value = "'" + '8.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.95%  '
$ws.Range("D21").Value = "'" + '371.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.71%  '
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").Value = "'" + '4.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.97%  '
$ws.Range("D24").Value = "'" + '71.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = "'" + '3.066.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.54%  '
$ws.Range("E27").Value = '  +1.49%  '
$ws.Range("D28").Value = "'" + '9.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("E29").Value = '  +6.88%  '
$ws.Range("D30").Value = "'" + '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("E31").Value = '  +1.79%  '
$ws.Range("D32").Value = "'" + '503.79'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.52%  '
$ws.Range("D33").Value = "'" + '7.74'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.29%  '
$ws.Range("E34").Value = '  +1.73%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = "'" + '165.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.51%  '
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("D38").Value = "'" + '19.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.64%  '
$ws.Range("E39").Value = '  -3.98%  '
$ws.Range("E40").Value = '  +20.26%  '
$ws.Range("D42").Value = "'" + '181.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.61%  '
$ws.Range("E43").Value = '  +3.05%  '
$ws.Range("D44").Value = "'" + '4.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").Value = "'" + '40.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("D47").Value = "'" + '1.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.27%  '
$ws.Range("D48").Value = "'" + '2.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").Value = "'" + '0.574'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.26%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("E51").Value = '  +3.28%  '
